$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 152
$ws.Range("A2").Value = 144
$ws.Range("A3").Value = 168
$ws.Range("A4").Value = 138
$ws.Range("A5").Value = 126
$ws.Range("A6").Value = 190
$ws.Range("A7").Value = 152.3999999999996
$ws.Range("A8").Value = 185.2000000000007
